# Insert two new weekly price records for "Brócoli" at the top of the
# date-sorted block (rows 324-387), pushing the existing records down by
# two rows (to 326-389). The two newly inserted rows capture a new
# observation dated 2023-01-25 (serial 44951) for "Primera" and "Segunda"
# quality grades.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 324 (existing row 324 and everything below
# shifts down by two rows).
$ws.Rows.Item(324).Insert()
$ws.Rows.Item(324).Insert()

# New row 324: Primera
$ws.Range("A324").Value = 7
$ws.Range("B324").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C324").Value = "Ñuble"
$ws.Range("D324").Value = 44951
$ws.Range("E324").Value = 16
$ws.Range("F324").Value = 100112023
$ws.Range("G324").Value = "Brócoli"
$ws.Range("H324").Value = "Sin especificar"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 300
$ws.Range("K324").Value = 1000
$ws.Range("L324").Value = 1000
$ws.Range("M324").Value = 1000
$ws.Range("N324").Value = "$/unidad"
$ws.Range("O324").Value = "Región del Maule"
$ws.Range("P324").Value = 1000
$ws.Range("Q324").Value = 1
$ws.Range("R324").Value = "Hortaliza"

# New row 325: Segunda
$ws.Range("A325").Value = 7
$ws.Range("B325").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C325").Value = "Ñuble"
$ws.Range("D325").Value = 44951
$ws.Range("E325").Value = 16
$ws.Range("F325").Value = 100112023
$ws.Range("G325").Value = "Brócoli"
$ws.Range("H325").Value = "Sin especificar"
$ws.Range("I325").Value = "Segunda"
$ws.Range("J325").Value = 400
$ws.Range("K325").Value = 700
$ws.Range("L325").Value = 750
$ws.Range("M325").Value = 725
$ws.Range("N325").Value = "$/unidad"
$ws.Range("O325").Value = "Región del Maule"
$ws.Range("P325").Value = 725
$ws.Range("Q325").Value = 1
$ws.Range("R325").Value = "Hortaliza"
